$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.139442682266235
$ws.Range("B1").Value = 4.760705471038818
$ws.Range("C1").Value = 3.319883346557617
$ws.Range("D1").Value = 2.241350650787354
$ws.Range("E1").Value = 2.055217742919922
